$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New summary rows appended below the per-instance result table (rows 2-11):
#   row 12  -> average of the "k" column (J)
#   row 14  -> Average of SW(S*)/SW(OPT)
#   row 15  -> Average of SC(S*)/SC(OPT)
#   row 16  -> Worst  of SW(S*)/SW(OPT)
#   row 17  -> Worst  of SC(S*)/SC(OPT)
# ---------------------------------------------------------------------------

# Row 14-17 labels + formulas (styled first so the bold/size-12/vertical-
# center style is created before the smaller bold style used by J12 -
# matches the style-table order produced by the original authoring session)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold / size-12 / vertically-centered style on a scratch cell and
# copy just the formatting onto the summary block - this avoids leaving any
# unused intermediate cell-format records behind in styles.xml.
$ws.Range("AB1").Font.Bold = $true
$ws.Range("AB1").Font.Size = 12
$ws.Range("AB1").VerticalAlignment = -4108
$ws.Range("AB1").Copy()
$ws.Range("A14:B17").PasteSpecial(-4122)
$ws.Range("AB1").Clear()

# Row 12: bold average of the k column
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

$ws.Range("AA1").Font.Bold = $true
$ws.Range("AA1").Copy()
$ws.Range("J12").PasteSpecial(-4122)
$ws.Range("AA1").Clear()

# Match the page setup seen on the authored copy (portrait / paper size 9 = A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the same cell selected as in the authored workbook
$ws.Range("J12").Select()
